$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '58.750.09'
$ws.Range("E2").Value = '  -6.53%  '
$ws.Range("D3").Value = '2.437.29'
$ws.Range("E3").Value = '  -9.16%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''535.65'
$ws.Range("E5").Value = '  -3.64%  '
$ws.Range("D6").Value = '''145.87'
$ws.Range("E6").Value = '  -7.54%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  -2.58%  '
$ws.Range("D9").Value = '2.448.01'
$ws.Range("E9").Value = '  -8.95%  '
$ws.Range("D10").Value = '''0.0987'
$ws.Range("E10").Value = '  -6.88%  '
$ws.Range("E11").Value = '  -2.15%  '
$ws.Range("D12").Value = '''5.30'
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("E13").Value = '  -4.98%  '
$ws.Range("D14").Value = '2.874.27'
$ws.Range("E14").Value = '  -8.95%  '
$ws.Range("D15").Value = '''23.81'
$ws.Range("E15").Value = '  -10.17%  '
$ws.Range("D16").Value = '58.664.88'
$ws.Range("E16").Value = '  -6.54%  '
$ws.Range("E17").Value = '  -6.18%  '
$ws.Range("D18").Value = '2.497.03'
$ws.Range("E18").Value = '  -7.03%  '
$ws.Range("D19").Value = '''11.09'
$ws.Range("D20").Value = '''4.34'
$ws.Range("E20").Value = '  -5.98%  '
$ws.Range("D21").Value = '''323.41'
$ws.Range("E21").Value = '  -6.25%  '
$ws.Range("D22").Value = '''0.965'
$ws.Range("E22").Value = '  -3.42%  '
$ws.Range("D23").Value = '''5.67'
$ws.Range("E23").Value = '  -8.71%  '
$ws.Range("D24").Value = '''60.49'
$ws.Range("E24").Value = '  -4.22%  '
$ws.Range("D25").Value = '''0.449'
$ws.Range("E25").Value = '  -11.27%  '
$ws.Range("D26").Value = '''0.159'
$ws.Range("E26").Value = '  -5.81%  '
$ws.Range("E27").Value = '  -2.44%  '
$ws.Range("D28").Value = '''7.62'
$ws.Range("E28").Value = '  -6.66%  '
$ws.Range("E29").Value = '  -6.67%  '
$ws.Range("D30").Value = '0.0₃0766'
$ws.Range("E30").Value = '  -10.45%  '
$ws.Range("D31").Value = '''6.60'
$ws.Range("E31").Value = '  -8.79%  '
$ws.Range("D32").Value = '''1.20'
$ws.Range("E32").Value = '  -14.53%  '
$ws.Range("D34").Value = '''156.14'
$ws.Range("E34").Value = '  -4.87%  '
$ws.Range("D35").Value = '''18.47'
$ws.Range("E35").Value = '  -5.15%  '
$ws.Range("E36").Value = '  -7.32%  '
$ws.Range("D37").Value = '''4.40'
$ws.Range("E37").Value = '  -10.21%  '
$ws.Range("D38").Value = '''1.67'
$ws.Range("E38").Value = '  -5.99%  '
$ws.Range("D39").Value = '''5.78'
$ws.Range("E39").Value = '  -6.61%  '
$ws.Range("D40").Value = '''310.91'
$ws.Range("E40").Value = '  -10.45%  '
$ws.Range("D41").Value = '''36.14'
$ws.Range("D42").Value = '''0.828'
$ws.Range("E42").Value = '  -12.05%  '
$ws.Range("D43").Value = '''3.68'
$ws.Range("E43").Value = '  -7.55%  '
$ws.Range("D44").Value = '''0.997'
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").Value = '''10.74'
$ws.Range("E45").Value = '  -2.31%  '
$ws.Range("D46").Value = '''0.0939'
$ws.Range("E46").Value = '  -3.13%  '
$ws.Range("D47").Value = '''0.578'
$ws.Range("E47").Value = '  -6.47%  '
$ws.Range("D48").Value = '''0.0523'
$ws.Range("E48").Value = '  -6.13%  '
$ws.Range("E49").Value = '  -5.18%  '
$ws.Range("D50").Value = '''121.61'
$ws.Range("E50").Value = '  -5.50%  '
$ws.Range("D51").Value = '''18.34'
$ws.Range("E51").Value = '  -8.72%  '
